# This script applies the cryptos-list price/volume update described in the commit,
# dated "Fri Nov 10 08:26:23 UTC 2023", matching the author's automated GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: the "Price" column (D) sometimes holds purely numeric-looking
# text (e.g. "250.92") that Excel would otherwise auto-convert to a real number,
# losing formatting (e.g. trailing zero in "14.40"). Force the cell to Text format
# first, assign the literal string, then restore the default "Normal" style so no
# extra formatting is left behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Rows with only Price (D) and Volume(1h) (E) changes ---
Set-TextValue "D2" "36.474.00"
$ws.Range("E2").Value = "  -0.35%  "
Set-TextValue "D3" "2.095.96"
$ws.Range("E3").Value = "  +9.58%  "
Set-TextValue "D5" "250.92"
$ws.Range("E5").Value = "  +0.89%  "
Set-TextValue "D6" "0.655"
$ws.Range("E6").Value = "  -6.15%  "
Set-TextValue "D8" "47.25"
$ws.Range("E8").Value = "  +6.92%  "
Set-TextValue "D9" "59.27"
$ws.Range("E9").Value = "  +2.53%  "
Set-TextValue "D10" "0.373"
$ws.Range("E10").Value = "  +0.71%  "
Set-TextValue "D11" "0.0741"
$ws.Range("E11").Value = "  -2.64%  "
Set-TextValue "D12" "0.0998"
$ws.Range("E12").Value = "  -0.31%  "
Set-TextValue "D13" "2.400.95"
$ws.Range("E13").Value = "  +9.52%  "
Set-TextValue "D14" "14.40"
$ws.Range("E14").Value = "  -0.36%  "
Set-TextValue "D15" "0.822"
$ws.Range("E15").Value = "  +1.78%  "
Set-TextValue "D16" "2.095.11"
$ws.Range("E16").Value = "  +9.53%  "
Set-TextValue "D17" "5.07"
$ws.Range("E17").Value = "  -0.62%  "
Set-TextValue "D18" "36.430.68"
$ws.Range("E18").Value = "  -0.51%  "
Set-TextValue "D19" "72.46"
$ws.Range("E19").Value = "  -2.60%  "
Set-TextValue "D20" "0.0₃0825"
$ws.Range("E20").Value = "  -4.01%  "
Set-TextValue "D21" "13.19"
$ws.Range("E21").Value = "  -1.18%  "
Set-TextValue "D22" "238.98"
$ws.Range("E22").Value = "  -4.36%  "
Set-TextValue "D23" "5.12"
$ws.Range("E23").Value = "  -1.38%  "
Set-TextValue "D26" "170.17"
$ws.Range("E26").Value = "  +0.93%  "
Set-TextValue "D27" "21.28"
$ws.Range("E27").Value = "  +13.91%  "
Set-TextValue "D28" "9.07"
$ws.Range("E28").Value = "  +3.02%  "
Set-TextValue "D29" "1.96"
$ws.Range("E29").Value = "  -10.15%  "
Set-TextValue "D30" "28.29"
$ws.Range("E30").Value = "  +57.54%  "
Set-TextValue "D31" "0.122"
$ws.Range("E31").Value = "  -4.99%  "
Set-TextValue "D34" "0.0924"
$ws.Range("E34").Value = "  +2.82%  "
Set-TextValue "D35" "0.948"
$ws.Range("E35").Value = "  +8.15%  "
Set-TextValue "D39" "4.05"
$ws.Range("E39").Value = "  -6.56%  "
Set-TextValue "D41" "1.16"
$ws.Range("E41").Value = "  +5.31%  "
Set-TextValue "D42" "0.0221"
$ws.Range("E42").Value = "  -2.48%  "
Set-TextValue "D43" "97.29"
$ws.Range("E43").Value = "  -8.64%  "
Set-TextValue "D44" "2.75"
$ws.Range("E44").Value = "  -6.87%  "
Set-TextValue "D45" "15.96"
$ws.Range("E45").Value = "  -8.11%  "
Set-TextValue "D46" "1.328.61"
$ws.Range("E46").Value = "  -1.29%  "
Set-TextValue "D47" "0.0838"
$ws.Range("E47").Value = "  +3.22%  "
Set-TextValue "D48" "6.93"
$ws.Range("E48").Value = "  +8.64%  "
Set-TextValue "D51" "2.22"
$ws.Range("E51").Value = "  -6.70%  "

# --- Rows with only Volume(1h) (E) changes ---
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -5.97%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -12.30%  "

# --- Rows where two coins swapped ranking position (Coin, Link, Price, Volume all change) ---
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "4.42"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0611"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D37" "2.32"
$ws.Range("E37").Value = "  +14.70%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D38" "1.87"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D49" "2.83"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D50" "2.273.89"
$ws.Range("E50").Value = "  +8.60%  "
